$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 1202.0588
$ws.Cells.Item(107, 9).Value = 1204.7273
$ws.Cells.Item(107, 10).Value = 1197.1666
$ws.Cells.Item(107, 11).Value = 1204.7273
$ws.Cells.Item(107, 12).Value = 1197.1666
$ws.Cells.Item(107, 13).Value = 715.2727
$ws.Cells.Item(107, 14).Value = -5037.1666

$ws.Cells.Item(116, 8).Value = 6917.1665
$ws.Cells.Item(116, 9).Value = 5216
$ws.Cells.Item(116, 10).Value = 9298.799999999999
$ws.Cells.Item(116, 11).Value = 5216
$ws.Cells.Item(116, 12).Value = 9298.799999999999
$ws.Cells.Item(116, 13).Value = -1774
$ws.Cells.Item(116, 14).Value = -16182.8

$ws.Cells.Item(129, 8).Value = 168614.25
$ws.Cells.Item(129, 9).Value = 168614.25
$ws.Cells.Item(129, 11).Value = 505842.75
$ws.Cells.Item(129, 13).Value = -500842.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 45560290
$ws.Cells.Item(86, 9).Value = 55679016
$ws.Cells.Item(86, 11).Value = 55679016
$ws.Cells.Item(86, 13).Value = -55677893

$ws.Cells.Item(89, 8).Value = 45560290
$ws.Cells.Item(89, 9).Value = 55679016
$ws.Cells.Item(89, 11).Value = 278395080
$ws.Cells.Item(89, 13).Value = -278389464

$ws.Cells.Item(105, 8).Value = 2939.5217
$ws.Cells.Item(105, 9).Value = 1779.4286
$ws.Cells.Item(105, 10).Value = 4744.1113
$ws.Cells.Item(105, 11).Value = 1779.4286
$ws.Cells.Item(105, 12).Value = 4744.1113
$ws.Cells.Item(105, 13).Value = -32.42859999999996
$ws.Cells.Item(105, 14).Value = -8238.1113

$ws.Cells.Item(133, 8).Value = 74999.336
$ws.Cells.Item(133, 10).Value = 74999.336
$ws.Cells.Item(133, 12).Value = 74999.336
$ws.Cells.Item(133, 14).Value = -85119.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3185.25
$ws.Cells.Item(31, 9).Value = 2451.2727
$ws.Cells.Item(31, 10).Value = 4800
$ws.Cells.Item(31, 11).Value = 2451.2727
$ws.Cells.Item(31, 12).Value = 4800
$ws.Cells.Item(31, 13).Value = -2156.2727
$ws.Cells.Item(31, 14).Value = -5390

$ws.Cells.Item(34, 8).Value = 3185.25
$ws.Cells.Item(34, 9).Value = 2451.2727
$ws.Cells.Item(34, 10).Value = 4800
$ws.Cells.Item(34, 11).Value = 2451.2727
$ws.Cells.Item(34, 12).Value = 4800
$ws.Cells.Item(34, 13).Value = -2249.2727
$ws.Cells.Item(34, 14).Value = -5204

$ws.Cells.Item(58, 8).Value = 3464.7273
$ws.Cells.Item(58, 9).Value = 2095.8
$ws.Cells.Item(58, 10).Value = 4605.5
$ws.Cells.Item(58, 11).Value = 2095.8
$ws.Cells.Item(58, 12).Value = 4605.5
$ws.Cells.Item(58, 13).Value = -1892.8
$ws.Cells.Item(58, 14).Value = -5011.5

$ws.Cells.Item(86, 8).Value = 17995.666
$ws.Cells.Item(86, 9).Value = 6500
$ws.Cells.Item(86, 10).Value = 21280.143
$ws.Cells.Item(86, 11).Value = 6500
$ws.Cells.Item(86, 12).Value = 21280.143
$ws.Cells.Item(86, 13).Value = -5377
$ws.Cells.Item(86, 14).Value = -23526.143

$ws.Cells.Item(89, 8).Value = 17995.666
$ws.Cells.Item(89, 9).Value = 6500
$ws.Cells.Item(89, 10).Value = 21280.143
$ws.Cells.Item(89, 11).Value = 32500
$ws.Cells.Item(89, 12).Value = 106400.715
$ws.Cells.Item(89, 13).Value = -26884
$ws.Cells.Item(89, 14).Value = -117632.715

$ws.Cells.Item(134, 8).Value = 2543.6667
$ws.Cells.Item(134, 9).Value = 1517.2142
$ws.Cells.Item(134, 10).Value = 16914
$ws.Cells.Item(134, 11).Value = 4551.642599999999
$ws.Cells.Item(134, 12).Value = 50742
$ws.Cells.Item(134, 13).Value = -2016.642599999999
$ws.Cells.Item(134, 14).Value = -55812

$ws.Cells.Item(136, 8).Value = 3464.7273
$ws.Cells.Item(136, 9).Value = 2095.8
$ws.Cells.Item(136, 10).Value = 4605.5
$ws.Cells.Item(136, 11).Value = 6287.400000000001
$ws.Cells.Item(136, 12).Value = 13816.5
$ws.Cells.Item(136, 13).Value = -3737.400000000001
$ws.Cells.Item(136, 14).Value = -18916.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(81, 8).Value = 966.3333
$ws.Cells.Item(81, 9).Value = 699
$ws.Cells.Item(81, 10).Value = 1100
$ws.Cells.Item(81, 11).Value = 2097
$ws.Cells.Item(81, 12).Value = 3300
$ws.Cells.Item(81, 13).Value = -974
$ws.Cells.Item(81, 14).Value = -5546

$ws.Cells.Item(84, 8).Value = 966.3333
$ws.Cells.Item(84, 9).Value = 699
$ws.Cells.Item(84, 10).Value = 1100
$ws.Cells.Item(84, 11).Value = 6291
$ws.Cells.Item(84, 12).Value = 9900
$ws.Cells.Item(84, 13).Value = -675
$ws.Cells.Item(84, 14).Value = -21132

$ws.Cells.Item(98, 8).Value = 1231.7142
$ws.Cells.Item(98, 10).Value = 1016.9091
$ws.Cells.Item(98, 12).Value = 3050.7273
$ws.Cells.Item(98, 14).Value = -6046.7273

$ws.Cells.Item(133, 8).Value = 10000
$ws.Cells.Item(133, 10).Value = 10000
$ws.Cells.Item(133, 12).Value = 30000
$ws.Cells.Item(133, 14).Value = -40120

$ws.Cells.Item(134, 8).Value = 3700.2727
$ws.Cells.Item(134, 9).Value = 2170.3
$ws.Cells.Item(134, 11).Value = 6510.900000000001
$ws.Cells.Item(134, 13).Value = -1440.900000000001

$ws.Cells.Item(136, 8).Value = 2620.6155
$ws.Cells.Item(136, 9).Value = 1152.5714
$ws.Cells.Item(136, 11).Value = 3457.7142
$ws.Cells.Item(136, 13).Value = 1642.2858

$ws.Cells.Item(138, 8).Value = 7739.5
$ws.Cells.Item(138, 9).Value = 6999.8
$ws.Cells.Item(138, 10).Value = 7986.067
$ws.Cells.Item(138, 11).Value = 20999.4
$ws.Cells.Item(138, 12).Value = 23958.201
$ws.Cells.Item(138, 13).Value = -15859.4
$ws.Cells.Item(138, 14).Value = -34238.201

$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 13).ClearContents()

$ws.Cells.Item(140, 8).Value = 2107.6155
$ws.Cells.Item(140, 9).Value = 1881.7273
$ws.Cells.Item(140, 11).Value = 5645.1819
$ws.Cells.Item(140, 13).Value = -465.1818999999996

$ws.Cells.Item(141, 8).Value = 2579.8
$ws.Cells.Item(141, 9).Value = 2579.8
$ws.Cells.Item(141, 11).Value = 7739.400000000001
$ws.Cells.Item(141, 13).Value = -2559.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 14912.714
$ws.Cells.Item(97, 10).Value = 25497.5
$ws.Cells.Item(97, 12).Value = 25497.5
$ws.Cells.Item(97, 14).Value = -26489.5

$ws.Cells.Item(113, 8).Value = 10048.875
$ws.Cells.Item(113, 9).Value = 4599
$ws.Cells.Item(113, 10).Value = 15498.75
$ws.Cells.Item(113, 11).Value = 4599
$ws.Cells.Item(113, 12).Value = 15498.75
$ws.Cells.Item(113, 13).Value = -2429
$ws.Cells.Item(113, 14).Value = -19838.75

$ws.Cells.Item(122, 8).Value = 5740.2856
$ws.Cells.Item(122, 9).Value = 5613.6665
$ws.Cells.Item(122, 10).Value = 6500
$ws.Cells.Item(122, 11).Value = 16840.9995
$ws.Cells.Item(122, 12).Value = 19500
$ws.Cells.Item(122, 13).Value = -14390.9995
$ws.Cells.Item(122, 14).Value = -24400

$ws.Cells.Item(126, 8).Value = 3255
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 13).ClearContents()

$ws.Cells.Item(132, 8).Value = 7555.643
$ws.Cells.Item(132, 9).Value = 6916.303
$ws.Cells.Item(132, 10).Value = 9899.888999999999
$ws.Cells.Item(132, 11).Value = 20748.909
$ws.Cells.Item(132, 12).Value = 29699.667
$ws.Cells.Item(132, 13).Value = -18218.909
$ws.Cells.Item(132, 14).Value = -34759.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3127.0667
$ws.Cells.Item(46, 10).Value = 4007.6
$ws.Cells.Item(46, 12).Value = 4007.6
$ws.Cells.Item(46, 14).Value = -4383.6

$ws.Cells.Item(59, 8).Value = 0
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 12).Value = 0
$ws.Cells.Item(59, 14).ClearContents()

$ws.Cells.Item(62, 8).Value = 52178.43
$ws.Cells.Item(62, 10).Value = 52178.43
$ws.Cells.Item(62, 12).Value = 52178.43
$ws.Cells.Item(62, 14).Value = -53426.43

$ws.Cells.Item(65, 8).Value = 52178.43
$ws.Cells.Item(65, 10).Value = 52178.43
$ws.Cells.Item(65, 12).Value = 156535.29
$ws.Cells.Item(65, 14).Value = -162775.29

$ws.Cells.Item(76, 8).Value = 6261
$ws.Cells.Item(76, 9).Value = 6261
$ws.Cells.Item(76, 11).Value = 6261
$ws.Cells.Item(76, 13).Value = -5923

$ws.Cells.Item(79, 8).Value = 6261
$ws.Cells.Item(79, 9).Value = 6261
$ws.Cells.Item(79, 11).Value = 6261
$ws.Cells.Item(79, 13).Value = -5091

$ws.Cells.Item(94, 8).Value = 60000
$ws.Cells.Item(94, 10).Value = 60000
$ws.Cells.Item(94, 12).Value = 60000
$ws.Cells.Item(94, 14).Value = -61352

$ws.Cells.Item(109, 8).Value = 17784.688
$ws.Cells.Item(109, 10).Value = 17784.688
$ws.Cells.Item(109, 12).Value = 17784.688
$ws.Cells.Item(109, 14).Value = -20558.688

$ws.Cells.Item(127, 8).Value = 53545
$ws.Cells.Item(127, 10).Value = 53545
$ws.Cells.Item(127, 12).Value = 53545
$ws.Cells.Item(127, 14).Value = -63465

$ws.Cells.Item(132, 8).Value = 2725.6875
$ws.Cells.Item(132, 9).Value = 2497.2144
$ws.Cells.Item(132, 10).Value = 4325
$ws.Cells.Item(132, 11).Value = 7491.6432
$ws.Cells.Item(132, 12).Value = 12975
$ws.Cells.Item(132, 13).Value = -4961.6432
$ws.Cells.Item(132, 14).Value = -18035

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2018.4
$ws.Cells.Item(122, 9).Value = 2018.4
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 6055.200000000001
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -3605.200000000001
$ws.Cells.Item(122, 14).ClearContents()

$ws.Cells.Item(123, 8).Value = 39733
$ws.Cells.Item(123, 10).Value = 39733
$ws.Cells.Item(123, 12).Value = 39733
$ws.Cells.Item(123, 14).Value = -49533

$ws.Cells.Item(126, 8).Value = 1025.6666
$ws.Cells.Item(126, 9).Value = 989.5
$ws.Cells.Item(126, 10).Value = 1098
$ws.Cells.Item(126, 11).Value = 2968.5
$ws.Cells.Item(126, 12).Value = 3294
$ws.Cells.Item(126, 13).Value = -498.5
$ws.Cells.Item(126, 14).Value = -8234

$ws.Cells.Item(129, 8).Value = 68998.5
$ws.Cells.Item(129, 9).Value = 45000
$ws.Cells.Item(129, 10).Value = 76998
$ws.Cells.Item(129, 11).Value = 45000
$ws.Cells.Item(129, 12).Value = 76998
$ws.Cells.Item(129, 13).Value = -40000
$ws.Cells.Item(129, 14).Value = -86998

$ws.Cells.Item(132, 8).Value = 6006.1333
$ws.Cells.Item(132, 9).Value = 2863.7144
$ws.Cells.Item(132, 10).Value = 50000
$ws.Cells.Item(132, 11).Value = 8591.143199999999
$ws.Cells.Item(132, 12).Value = 150000
$ws.Cells.Item(132, 13).Value = -6061.143199999999
$ws.Cells.Item(132, 14).Value = -155060

$ws.Cells.Item(136, 8).Value = 89092.55
$ws.Cells.Item(136, 9).Value = 97251.8
$ws.Cells.Item(136, 10).Value = 7500
$ws.Cells.Item(136, 11).Value = 291755.4
$ws.Cells.Item(136, 12).Value = 22500
$ws.Cells.Item(136, 13).Value = -289205.4
$ws.Cells.Item(136, 14).Value = -27600
